$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13 (pushes old rows 13-24 down to 14-25,
# carrying their row heights / cell styles with them automatically).
$ws.Rows("13:13").Insert()

# The insert leaves a stray formatted/empty A13 cell (inherited from the row
# above) plus B13/C13 with the wrong inherited style; clear the row fully so we
# can rebuild it from scratch with the formatting copied from row 14 (B/C only,
# no A cell in the final layout for this row).
$ws.Range("A13:C13").Clear()

$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New row 13 content: professor name relocated here from old row 18 (B/C).
$ws.Range("B13").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C13").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# Content updates on the rows whose text changed (row numbers are already the
# POST-insert / final numbering).
$ws.Range("B10").Value = "Proporcionar aos alunos os fundamentos teóricos do tratamento e da destinação final dos resíduos sólidos produzidos em instituições, indústrias, comunidades ou municípios, contribuindo para a melhoria das condições ambientais nas cidades e no campo, bem como atualizar o conhecimento teórico de profissionais da área de gerenciamento de resíduos sólidos."
$ws.Range("C10").Value = "Proporcionar aos alunos os fundamentos teóricos do tratamento e da destinação final dos resíduos sólidos produzidos em instituições, indústrias, comunidades ou municípios, contribuindo para a melhoria das condições ambientais nas cidades e no campo, bem como atualizar o conhecimento teórico de profissionais da área de gerenciamento de resíduos sólidos."

$ws.Range("B14").Value = "Tipos de resíduos; Gestão e legislação no Brasil"
$ws.Range("C14").Value = "Tipos de resíduos; Gestão e legislação no Brasil"

$ws.Range("B16").Value = "Conceitos e definições de resíduos sólidos; amostragem, caracterização e classificação de resíduos sólidos; sistema de coleta, acondicionamento, transporte, tratamento e disposição final de resíduos sólidos; a Política Nacional de Resíduos Sólidos; modelo de gerenciamento de resíduos sólidos."
$ws.Range("C16").Value = "Conceitos e definições de resíduos sólidos; amostragem, caracterização e classificação de resíduos sólidos; sistema de coleta, acondicionamento, transporte, tratamento e disposição final de resíduos sólidos; a Política Nacional de Resíduos Sólidos; modelo de gerenciamento de resíduos sólidos."

$ws.Range("B19").Value = "As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático."
$ws.Range("C19").Value = "As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático."

$ws.Range("B20").Value = "Serão aplicadas duas avaliações para compor a média que será a soma das duas provas, sendo o resultado dividido por dois."
$ws.Range("C20").Value = "Serão aplicadas duas avaliações para compor a média que será a soma das duas provas, sendo o resultado dividido por dois."

$ws.Range("B21").Value = "A Nota Final será composta pela Média obtida da Nota do Período somada à Nota de Recuperação e dividido por dois"
$ws.Range("C21").Value = "A Nota Final será composta pela Média obtida da Nota do Período somada à Nota de Recuperação e dividido por dois"

$ws.Range("B22").Value = "Bibliografia básica:BARROS, R.M., Tratado sobre resíduos sólidos: gestão, uso e sustentabilidade, Editora Interciência, 2013.DA SILVA-FILHO, C.R.V., SOLER, F.D., Gestão de resíduos sólidos: o que diz a lei, 2° ed., Editora Trevisan, 2013.RIBEIRO, D.V., MORELLI, M.R., Resíduos sólidos: problemas ou oportunidades?, Editora Interciência, 2009.Bibliografia complementar:BRAGA B. (Org.), Introdução à engenharia ambiental: o desafio do desenvolvimento sustentável, 2° ed., Ed. Pearson Prentice Hall, 2005.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão, Ed. Campus, 2013.CARVALHO, T.C.M.B., XAVIER, L.H. (Org.), Gestão de resíduos eletroeletrônicos: uma abordagem prática para a sustentabilidade, Edidora Elsivier Ltda, 2014.JACOBI, P. (Org.), Gestão compartilhada dos resíduos sólidos no Brasil: inovação com inclusão social, Annablume, 2006PEREIRA-NETO, J.T., Gerenciamento do lixo urbano: aspectos técnicos e operacionais, Editora UFV, 2013.SÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos, Ed. Oficina de textos, 3° reimpressão, 2011ASSOCIAÇÃO BRASILEIRA DE NORMAS TÉCNICAS (ABNT) – NORMAS ABNT RESÍDUOS SÓLIDOS: COLETÂNEA DE NORMAS: NBR 10.004, NBR 10.005, NBR 10.006, NBR 10.007ASSOCIAÇÃO BRSILEIRA DE NORMAS TÉCNICAS (ABNT) – NORMAS TÉCNICAS APRESENTAÇÃO DE PROJETOS DE ATERROS CONTROLADOS DE RESÍDUOS SÓLIDOS URBANOS: NBR 8849"
$ws.Range("C22").Value = "Bibliografia básica:BARROS, R.M., Tratado sobre resíduos sólidos: gestão, uso e sustentabilidade, Editora Interciência, 2013.DA SILVA-FILHO, C.R.V., SOLER, F.D., Gestão de resíduos sólidos: o que diz a lei, 2° ed., Editora Trevisan, 2013.RIBEIRO, D.V., MORELLI, M.R., Resíduos sólidos: problemas ou oportunidades?, Editora Interciência, 2009.Bibliografia complementar:BRAGA B. (Org.), Introdução à engenharia ambiental: o desafio do desenvolvimento sustentável, 2° ed., Ed. Pearson Prentice Hall, 2005.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão, Ed. Campus, 2013.CARVALHO, T.C.M.B., XAVIER, L.H. (Org.), Gestão de resíduos eletroeletrônicos: uma abordagem prática para a sustentabilidade, Edidora Elsivier Ltda, 2014.JACOBI, P. (Org.), Gestão compartilhada dos resíduos sólidos no Brasil: inovação com inclusão social, Annablume, 2006PEREIRA-NETO, J.T., Gerenciamento do lixo urbano: aspectos técnicos e operacionais, Editora UFV, 2013.SÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos, Ed. Oficina de textos, 3° reimpressão, 2011ASSOCIAÇÃO BRASILEIRA DE NORMAS TÉCNICAS (ABNT) – NORMAS ABNT RESÍDUOS SÓLIDOS: COLETÂNEA DE NORMAS: NBR 10.004, NBR 10.005, NBR 10.006, NBR 10.007ASSOCIAÇÃO BRSILEIRA DE NORMAS TÉCNICAS (ABNT) – NORMAS TÉCNICAS APRESENTAÇÃO DE PROJETOS DE ATERROS CONTROLADOS DE RESÍDUOS SÓLIDOS URBANOS: NBR 8849"

